$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update price (D) and volume-change (E) columns,
# plus the row-34/35 coin swap (EthereumClassic <-> NEARProtocol).

$ws.Range('D2').Value = '63.817.81'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.751.78'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '580.33'
$ws.Range('E5').Value = '  -2.50%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '158.22'
$ws.Range('E6').Value = '  +4.04%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.96'
$ws.Range('E10').Value = '  -12.12%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.391'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '3.239.71'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.87'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').Value = '63.795.33'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000155'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = '2.763.85'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.24'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.95'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '361.30'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.88'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.566'
$ws.Range('E22').Value = '  +4.64%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.992'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.48'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.67'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = '0.0₃0928'
$ws.Range('E28').Value = '  +4.92%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.00'
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.11'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('E31').Value = '  +3.23%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '168.85'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '20.55'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.01'
$ws.Range('E35').Value = '  +4.67%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.46'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.81'
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.20'
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.14'
$ws.Range('E40').Value = '  +7.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '330.71'
$ws.Range('E41').Value = '  -6.67%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.66'
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '22.04'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0599'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.87'
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0257'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '136.06'
$ws.Range('E48').Value = '  -6.99%  '
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.05'
$ws.Range('E51').Value = '  +0.95%  '
